# Eggerland - several more levels done.
# Append rows 47-54 to Sheet1 (mirrors the existing Enter/Exit-room timing log)
# and extend the D-column "duration" formula down through the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @(47, "Get Key",        82883, 81079),
    @(48, "Get Key",        84163, 82358),
    @(49, "Get Key",        85727, 83923),
    @(50, "Get heart",      86849, 85045),
    @(51, "Get key",        87329, 85525),
    @(52, "Square appears", 89980, 88177),
    @(53, "Get key",        90538, 88736),
    @(54, "Last heart",     91107, 89305)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Fill column D (=C-B) for the newly added rows as one shared-formula block.
$ws.Range("D47:D54").Formula = "=C47-B47"

# Update the view to match the post-edit state (bottom of the new data).
$ws.Range("B55").Select() | Out-Null
